# Auto-generated script to apply cryptos list price/volume update
# (GitHub Actions scheduled data refresh simulation)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The Price column (D) holds values that look numeric (e.g. '1.003') but
# must remain plain text, matching the original inlineStr cell type.
# Temporarily force Text format on the whole Price column so Excel does not
# reinterpret these assignments as numbers, then restore the original
# (default/Normal) style so no stray formatting is introduced.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range('D2').Value = '28.853.17'
$ws.Range('E2').Value = '  -0.39%  '
$ws.Range('D3').Value = '1.917.26'
$ws.Range('E3').Value = '  +0.72%  '
$ws.Range('D4').Value = '1.003'
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').Value = '324.34'
$ws.Range('E5').Value = '  +0.09%  '
$ws.Range('D6').Value = '1.002'
$ws.Range('E6').Value = '  -0.03%  '
$ws.Range('D7').Value = '0.4562'
$ws.Range('E7').Value = '  -0.72%  '
$ws.Range('D8').Value = '0.3806'
$ws.Range('E8').Value = '  -0.15%  '
$ws.Range('D9').Value = '0.07766'
$ws.Range('E9').Value = '  +0.83%  '
$ws.Range('D10').Value = '0.9775'
$ws.Range('E10').Value = '  -0.35%  '
$ws.Range('D11').Value = '22.24'
$ws.Range('E11').Value = '  +0.96%  '
$ws.Range('D12').Value = '1.944.35'
$ws.Range('E12').Value = '  +2.29%  '
$ws.Range('D13').Value = '5.700'
$ws.Range('E13').Value = '  +0.47%  '
$ws.Range('D14').Value = '6.977'
$ws.Range('E14').Value = '  +0.33%  '
$ws.Range('D15').Value = '0.06988'
$ws.Range('E15').Value = '  -0.77%  '
$ws.Range('D16').Value = '1.004'
$ws.Range('E16').Value = '  -0.03%  '
$ws.Range('D17').Value = '84.42'
$ws.Range('E17').Value = '  +0.79%  '
$ws.Range('D18').Value = '0.000009495'
$ws.Range('E18').Value = '  -0.31%  '
$ws.Range('D19').Value = '16.67'
$ws.Range('E19').Value = '  -0.22%  '
$ws.Range('D20').Value = '1.002'
$ws.Range('E20').Value = '  +0.00%  '
$ws.Range('D21').Value = '28.867.38'
$ws.Range('E21').Value = '  -0.27%  '
$ws.Range('D22').Value = '5.338'
$ws.Range('E22').Value = '  +0.31%  '
$ws.Range('D23').Value = '11.09'
$ws.Range('E23').Value = '  +1.86%  '
$ws.Range('D24').Value = '2.149.41'
$ws.Range('E24').Value = '  +0.88%  '
$ws.Range('D25').Value = '2.058'
$ws.Range('E25').Value = '  -1.87%  '
$ws.Range('D26').Value = '157.90'
$ws.Range('E26').Value = '  +0.59%  '
$ws.Range('D27').Value = '19.06'
$ws.Range('E27').Value = '  +0.04%  '
$ws.Range('D28').Value = '5.626'
$ws.Range('E28').Value = '  +0.52%  '
$ws.Range('D29').Value = '117.68'
$ws.Range('E29').Value = '  +0.10%  '
$ws.Range('D31').Value = '0.09288'
$ws.Range('E31').Value = '  +0.04%  '
$ws.Range('D32').Value = '0.8716'
$ws.Range('E32').Value = '  +1.12%  '
$ws.Range('D33').Value = '5.105'
$ws.Range('E33').Value = '  +0.64%  '
$ws.Range('D34').Value = '1.244'
$ws.Range('E34').Value = '  -0.38%  '
$ws.Range('D35').Value = '3.027'
$ws.Range('E35').Value = '  +0.71%  '
$ws.Range('D36').Value = '0.05704'
$ws.Range('E36').Value = '  +0.19%  '
$ws.Range('D39').Value = '0.02036'
$ws.Range('E39').Value = '  -0.01%  '
$ws.Range('D40').Value = '3.062'
$ws.Range('E40').Value = '  +11.60%  '
$ws.Range('D41').Value = '7.497'
$ws.Range('E41').Value = '  +0.53%  '
$ws.Range('D42').Value = '0.5506'
$ws.Range('E42').Value = '  -0.09%  '
$ws.Range('D43').Value = '0.1758'
$ws.Range('E43').Value = '  +0.03%  '
$ws.Range('D46').Value = '2.173'
$ws.Range('E46').Value = '  +3.86%  '
$ws.Range('D47').Value = '0.5156'
$ws.Range('E47').Value = '  -0.71%  '
$ws.Range('D48').Value = '0.06925'
$ws.Range('E48').Value = '  +1.68%  '
$ws.Range('D49').Value = '11.18'
$ws.Range('E49').Value = '  -1.44%  '
$ws.Range('D50').Value = '110.63'
$ws.Range('E50').Value = '  -0.64%  '
$ws.Range('D51').Value = '1.767'
$ws.Range('E51').Value = '  -0.70%  '
$ws.Range('E30').Value = '  +0.00%  '
$ws.Range('E37').Value = '  +0.10%  '
$ws.Range('E38').Value = '  +0.03%  '

# Row 44/45 swap: PEPE moves to row 44 (with updated data), Aptos moves to row 45 (with updated data)
$ws.Range("B44").Value = 'PEPE'
$ws.Range("C44").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D44").Value = '0.000002896'
$ws.Range("E44").Value = '  +16.96%  '
$ws.Range("B45").Value = 'Aptos'
$ws.Range("C45").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D45").Value = '9.340'
$ws.Range("E45").Value = '  +0.77%  '

# Restore the Price column's original (unset/Normal) style now that the
# text values are locked in, so no extra formatting/style diffs remain.
$priceRange.Style = "Normal"
